# Auto-generated script applying market-price refresh updates to all 8 sheets
# (columns H-N: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ))
$wb = $excel.ActiveWorkbook


# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 10209
$ws.Range("I11").Value = 10209
$ws.Range("K11").Value = 10209
$ws.Range("M11").Value = -10069
$ws.Range("H18").Value = 3505.4443
$ws.Range("I18").Value = 3505.4443
$ws.Range("K18").Value = 3505.4443
$ws.Range("M18").Value = -3221.4443
$ws.Range("H53").Value = 66667016
$ws.Range("I53").Value = 64
$ws.Range("K53").Value = 64
$ws.Range("M53").Value = 573
$ws.Range("H74").Value = 5944.222
$ws.Range("I74").Value = 5944.222
$ws.Range("K74").Value = 5944.222
$ws.Range("M74").Value = -5008.222
$ws.Range("H77").Value = 5944.222
$ws.Range("I77").Value = 5944.222
$ws.Range("K77").Value = 29721.11
$ws.Range("M77").Value = -25041.11
$ws.Range("H98").Value = 3290.611
$ws.Range("I98").Value = 2319.6667
$ws.Range("J98").Value = 5232.5
$ws.Range("K98").Value = 2319.6667
$ws.Range("L98").Value = 5232.5
$ws.Range("M98").Value = -821.6667000000002
$ws.Range("N98").Value = -8228.5
$ws.Range("H122").Value = 3290.611
$ws.Range("I122").Value = 2319.6667
$ws.Range("J122").Value = 5232.5
$ws.Range("K122").Value = 6959.000100000001
$ws.Range("L122").Value = 15697.5
$ws.Range("M122").Value = -4509.000100000001
$ws.Range("N122").Value = -20597.5
$ws.Range("H137").Value = 3711689.8
$ws.Range("I137").Value = 5875
$ws.Range("K137").Value = 17625
$ws.Range("M137").Value = -15075

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 2903957.8
$ws.Range("I61").Value = 67328.64999999999
$ws.Range("K61").Value = 67328.64999999999
$ws.Range("M61").Value = -67116.64999999999
$ws.Range("H74").Value = 637287.5
$ws.Range("I74").Value = 985.5
$ws.Range("J74").Value = 2025582.8
$ws.Range("K74").Value = 985.5
$ws.Range("L74").Value = 2025582.8
$ws.Range("M74").Value = -111.5
$ws.Range("N74").Value = -2027330.8
$ws.Range("H77").Value = 637287.5
$ws.Range("I77").Value = 985.5
$ws.Range("J77").Value = 2025582.8
$ws.Range("K77").Value = 4927.5
$ws.Range("L77").Value = 10127914
$ws.Range("M77").Value = -559.5
$ws.Range("N77").Value = -10136650
$ws.Range("H102").Value = 2354.1333
$ws.Range("I102").Value = 2354.1333
$ws.Range("K102").Value = 2354.1333
$ws.Range("M102").Value = -732.1333
$ws.Range("H122").Value = 2461.8462
$ws.Range("I122").Value = 2250.3333
$ws.Range("K122").Value = 6750.999899999999
$ws.Range("M122").Value = -4300.999899999999
$ws.Range("H123").Value = 143323.33
$ws.Range("J123").Value = 143323.33
$ws.Range("L123").Value = 143323.33
$ws.Range("N123").Value = -153123.33
$ws.Range("H136").Value = 2903957.8
$ws.Range("I136").Value = 67328.64999999999
$ws.Range("K136").Value = 201985.95
$ws.Range("M136").Value = -199435.95

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H22").Value = 1404.8
$ws.Range("I22").Value = 1255.3334
$ws.Range("K22").Value = 1255.3334
$ws.Range("M22").Value = -1082.3334
$ws.Range("H105").Value = 99998.5
$ws.Range("I105").Value = 99998.5
$ws.Range("K105").Value = 99998.5
$ws.Range("M105").Value = -98251.5
$ws.Range("H122").Value = 51998.7
$ws.Range("J122").Value = 51998.7
$ws.Range("L122").Value = 51998.7
$ws.Range("N122").Value = -61798.7
$ws.Range("H134").Value = 26472626
$ws.Range("I134").Value = 1617.1428
$ws.Range("K134").Value = 4851.428400000001
$ws.Range("M134").Value = -2316.428400000001

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 11905131
$ws.Range("I16").Value = 12987398
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 12987398
$ws.Range("L16").Value = 200
$ws.Range("M16").Value = -12987111
$ws.Range("N16").Value = -774
$ws.Range("H113").Value = 11905131
$ws.Range("I113").Value = 12987398
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 12987398
$ws.Range("L113").Value = 200
$ws.Range("M113").Value = -12985228
$ws.Range("N113").Value = -4540
$ws.Range("H122").Value = 3405.6155
$ws.Range("I122").Value = 3068.0908
$ws.Range("K122").Value = 9204.2724
$ws.Range("M122").Value = -6754.2724
$ws.Range("H134").Value = 1954.6
$ws.Range("I134").Value = 1491.3572
$ws.Range("J134").Value = 3035.5
$ws.Range("K134").Value = 4474.071599999999
$ws.Range("L134").Value = 9106.5
$ws.Range("M134").Value = -1939.071599999999
$ws.Range("N134").Value = -14176.5

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 1051.421
$ws.Range("I2").Value = 526.36365
$ws.Range("J2").Value = 1265.3334
$ws.Range("K2").Value = 3158.1819
$ws.Range("L2").Value = 7592.0004
$ws.Range("M2").Value = -3045.1819
$ws.Range("N2").Value = -7818.0004
$ws.Range("H4").Value = 1134419.1
$ws.Range("I4").Value = 940804.1
$ws.Range("J4").Value = 2167033
$ws.Range("K4").Value = 2822412.3
$ws.Range("L4").Value = 6501099
$ws.Range("M4").Value = -2822300.3
$ws.Range("N4").Value = -6501323
$ws.Range("H94").Value = 12552.2
$ws.Range("J94").Value = 17330.334
$ws.Range("L94").Value = 51991.00199999999
$ws.Range("N94").Value = -53343.00199999999
$ws.Range("H104").Value = 7356.5557
$ws.Range("J104").Value = 12330.25
$ws.Range("L104").Value = 36990.75
$ws.Range("N104").Value = -42232.75
$ws.Range("H107").Value = 1240.8125
$ws.Range("J107").Value = 1440.75
$ws.Range("L107").Value = 4322.25
$ws.Range("N107").Value = -8162.25
$ws.Range("H126").Value = 5000
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -24880
$ws.Range("H132").Value = 2443.1875
$ws.Range("I132").Value = 2409.1
$ws.Range("K132").Value = 21681.9
$ws.Range("M132").Value = -19151.9

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 100001720
$ws.Range("I102").Value = 125001624
$ws.Range("J102").Value = 2095
$ws.Range("K102").Value = 125001624
$ws.Range("L102").Value = 2095
$ws.Range("M102").Value = -125000002
$ws.Range("N102").Value = -5339
$ws.Range("H132").Value = 1070292.6
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1084.6818
$ws.Range("I16").Value = 1014.3
$ws.Range("J16").Value = 1788.5
$ws.Range("K16").Value = 1014.3
$ws.Range("L16").Value = 1788.5
$ws.Range("M16").Value = -844.3
$ws.Range("N16").Value = -2128.5
$ws.Range("H40").Value = 5126.5557
$ws.Range("I40").Value = 3810.2307
$ws.Range("J40").Value = 8549
$ws.Range("K40").Value = 3810.2307
$ws.Range("L40").Value = 8549
$ws.Range("M40").Value = -3674.2307
$ws.Range("N40").Value = -8821
$ws.Range("H59").Value = 160000
$ws.Range("J59").Value = 160000
$ws.Range("L59").Value = 160000
$ws.Range("N59").Value = -161308
$ws.Range("H76").Value = 15000
$ws.Range("I76").Value = 15000
$ws.Range("J76").Value = 15000
$ws.Range("K76").Value = 15000
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15676
$ws.Range("M76").Value = -14662
$ws.Range("H79").Value = 15000
$ws.Range("I79").Value = 15000
$ws.Range("J79").Value = 15000
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17340
$ws.Range("M79").Value = -13830
$ws.Range("H123").Value = 129990
$ws.Range("J123").Value = 129990
$ws.Range("L123").Value = 129990
$ws.Range("N123").Value = -139790
$ws.Range("H132").Value = 2723
$ws.Range("I132").Value = 2378.2778
$ws.Range("J132").Value = 3498.625
$ws.Range("K132").Value = 7134.8334
$ws.Range("L132").Value = 10495.875
$ws.Range("M132").Value = -4604.8334
$ws.Range("N132").Value = -15555.875

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H16").Value = 38200
$ws.Range("J16").Value = 38200
$ws.Range("L16").Value = 38200
$ws.Range("N16").Value = -38784
$ws.Range("H82").Value = 39999.5
$ws.Range("J82").Value = 39999.5
$ws.Range("L82").Value = 39999.5
$ws.Range("N82").Value = -40765.5
$ws.Range("H85").Value = 39999.5
$ws.Range("J85").Value = 39999.5
$ws.Range("L85").Value = 39999.5
$ws.Range("N85").Value = -42651.5
$ws.Range("H122").Value = 1915.2142
$ws.Range("I122").Value = 1676.1666
$ws.Range("J122").Value = 3349.5
$ws.Range("K122").Value = 5028.4998
$ws.Range("L122").Value = 10048.5
$ws.Range("M122").Value = -2578.4998
$ws.Range("N122").Value = -14948.5
$ws.Range("H132").Value = 2373.96
$ws.Range("I132").Value = 1658.579
$ws.Range("J132").Value = 4639.3335
$ws.Range("K132").Value = 4975.737
$ws.Range("L132").Value = 13918.0005
$ws.Range("M132").Value = -2445.737
$ws.Range("N132").Value = -18978.0005
$ws.Range("H141").Value = 115014.25
$ws.Range("J141").Value = 115014.25
$ws.Range("L141").Value = 115014.25
$ws.Range("N141").Value = -125374.25
